$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")

# New expense entry: Derek, 2023-10-13 (serial 45212), 302.47
$ws.Range("A14").Value = "Derek"
$ws.Range("B14").Value = 45212
$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing date style
$ws.Range("C14").Value = 302.47

# Make Budget the active sheet/tab with the given selection
$ws.Activate()
$null = $ws.Range("E21").Select()
